# Applies the attendance_reports sync edits to the "Session Analysis Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 - new recorder added + attendance count updated
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("H4").Value = "73/251"

# Row 10 - updated average attendance % (force as text, not an auto-converted percentage number,
# by entering it as a formula that evaluates to the text and then flattening to a plain value via
# paste-special so the original cell style/number format is preserved)
$ws.Range("L10").Formula = '="25.5%"'
$ws.Range("L10").Copy()
$ws.Range("L10").PasteSpecial(-4163)

# Row 15 - reorder "Recorded By" list + updated average attendance %
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("S15").Formula = '="25.5%"'
$ws.Range("S15").Copy()
$ws.Range("S15").PasteSpecial(-4163)

# Row 28 - reorder "Recorded By" list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

$excel.CutCopyMode = $false
